# Aggiornamento dati Concordia: inserita la riga mancante del 2021-02-08
# (44235), ricalcolata la somma mobile a 7 giorni (colonne C/D) per le righe
# vicine, e aggiunta la nuova riga del 2021-03-02 (44257) in fondo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) La riga 90 (05/02/2021) non si sposta, ma la somma mobile cambia.
$ws.Cells.Item(90, 3).Value = 7
$ws.Cells.Item(90, 4).Value = 85.07535245503161

# 2) Inserisce una nuova riga alla posizione 93 per la data mancante
#    (08/02/2021, seriale 44235), spostando in basso le righe successive.
$ws.Rows.Item(93).Insert()

# Copia il formato (stile data) della cella soprastante sulla nuova cella A93.
$ws.Cells.Item(92, 1).Copy()
$ws.Cells.Item(93, 1).PasteSpecial(-4122)

$ws.Cells.Item(93, 1).Value = 44235
$ws.Cells.Item(93, 2).Value = 0
$ws.Cells.Item(93, 3).Value = 7
$ws.Cells.Item(93, 4).Value = 85.07535245503161

# 3) Le righe 94-111 (ex 93-110) mantengono la stessa data/B ma la somma
#    mobile a 7 gg. (C/D) viene ricalcolata con il nuovo dato disponibile.
$ws.Cells.Item(94, 3).Value = 8
$ws.Cells.Item(94, 4).Value = 97.22897423432183

$ws.Cells.Item(95, 3).Value = 5
$ws.Cells.Item(95, 4).Value = 60.76810889645115

$ws.Cells.Item(96, 3).Value = 5
$ws.Cells.Item(96, 4).Value = 60.76810889645115

$ws.Cells.Item(97, 3).Value = 6
$ws.Cells.Item(97, 4).Value = 72.92173067574137

$ws.Cells.Item(98, 3).Value = 8
$ws.Cells.Item(98, 4).Value = 97.22897423432183

$ws.Cells.Item(99, 3).Value = 9
$ws.Cells.Item(99, 4).Value = 109.3825960136121

$ws.Cells.Item(100, 3).Value = 9
$ws.Cells.Item(100, 4).Value = 109.3825960136121

$ws.Cells.Item(101, 3).Value = 6
$ws.Cells.Item(101, 4).Value = 72.92173067574137

$ws.Cells.Item(102, 3).Value = 7
$ws.Cells.Item(102, 4).Value = 85.07535245503161

$ws.Cells.Item(103, 3).Value = 8
$ws.Cells.Item(103, 4).Value = 97.22897423432183

$ws.Cells.Item(104, 3).Value = 11
$ws.Cells.Item(104, 4).Value = 133.6898395721925

$ws.Cells.Item(105, 3).Value = 8
$ws.Cells.Item(105, 4).Value = 97.22897423432183

$ws.Cells.Item(106, 3).Value = 7
$ws.Cells.Item(106, 4).Value = 85.07535245503161

$ws.Cells.Item(107, 3).Value = 7
$ws.Cells.Item(107, 4).Value = 85.07535245503161

$ws.Cells.Item(108, 3).Value = 7
$ws.Cells.Item(108, 4).Value = 85.07535245503161

$ws.Cells.Item(109, 3).Value = 7
$ws.Cells.Item(109, 4).Value = 85.07535245503161

$ws.Cells.Item(110, 3).Value = 8
$ws.Cells.Item(110, 4).Value = 97.22897423432183

$ws.Cells.Item(111, 3).Value = 6
$ws.Cells.Item(111, 4).Value = 72.92173067574137

# 4) La riga 112 (ex 111, 27/02/2021) in precedenza non aveva somma mobile
#    calcolata: ora che ci sono abbastanza dati successivi viene valorizzata.
$ws.Cells.Item(112, 3).Value = 6
$ws.Cells.Item(112, 4).Value = 72.92173067574137

# 5) Aggiunge la nuova riga finale per il 02/03/2021 (seriale 44257).
$ws.Cells.Item(114, 1).Copy()
$ws.Cells.Item(115, 1).PasteSpecial(-4122)

$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 0

Write-Output "ok"
